$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 15: change year from 2046 to 2045
$ws.Range("D15").Value = 2045

# Row 16: populate like row 15/13 (new CAP_BND / LO bound row) but for 2045
# with the new process SUPDACELC5N driving the bound (forces DAC build)
$ws.Range("D16").Value = 2045
$ws.Range("E16").Value = "LO"
$ws.Range("F16").Value = "CAP_BND"
$ws.Range("G16").Value = $null
$ws.Range("H16").Value = $null
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1000
$ws.Range("M16").Value = "SUPDACELC5N"

# Match styling used by the rest of the data block (row 15) for D16:J16
$ws.Range("D16:J16").Style = $ws.Range("D15:J15").Style

# Update the active selection to D17 as recorded in the saved view state
$ws.Range("D17").Select()
